$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the descriptive header text cells (company name / report title / period)
$ws.Range("E2:E8").ClearContents()

# Clear the "mapped category" helper labels in column H (rows 11-54)
$ws.Range("H11:H54").ClearContents()

# Clear the checking formula in H56 (G56-F56)
$ws.Range("H56").ClearContents()

# Update the selected range shown when the file was last saved
$ws.Range("H1:H1048576").Select()
